# Generate Report for Archive
#
# 1. The "Status" value shown for this file moves from "Ready for handoff"
#    to "In Translation" everywhere it is displayed:
#      - Overview sheet: columns E (zh-cn) and F (de-de), row 2
#      - zh-cn sheet:      column C (Status), row 2
#      - de-de sheet:      column C (Status), row 2
# 2. Because the new status text is shorter, the Status column on each
#    sheet is re-sized (narrower) to fit the new content.

$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$zhcn     = $wb.Worksheets.Item("zh-cn")
$dede     = $wb.Worksheets.Item("de-de")

$newStatus = "In Translation"

# --- update the displayed status text -------------------------------------
$overview.Range("E2").Value = $newStatus
$overview.Range("F2").Value = $newStatus
$zhcn.Range("C2").Value = $newStatus
$dede.Range("C2").Value = $newStatus

# --- shrink the now-narrower Status columns to fit -------------------------
$overview.Columns.Item(5).ColumnWidth = 12.5   # column E (zh-cn status)
$overview.Columns.Item(6).ColumnWidth = 12.5   # column F (de-de status)
$zhcn.Columns.Item(3).ColumnWidth = 12.5        # column C (Status)
$dede.Columns.Item(3).ColumnWidth = 12.5        # column C (Status)
